$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain plain text so that numeric-looking
# strings like "0.999" are not silently converted into Excel numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Update Price (D) and Volume(1h) (E) cells row by row ---
$ws.Range("D2").Value = "67.985.45"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.336.07"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "583.73"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "177.81"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("E9").Value = "  +4.40%  "
$ws.Range("D10").Value = "0.583"
$ws.Range("D11").Value = "47.84"
$ws.Range("E11").Value = "  +5.38%  "
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("D13").Value = "698.11"
$ws.Range("E13").Value = "  +5.30%  "
$ws.Range("D14").Value = "3.875.31"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").Value = "8.47"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "68.013.71"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "3.336.28"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "11.16"
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("D21").Value = "0.896"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "5.40"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").Value = "16.98"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "100.48"
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("D25").Value = "3.92"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "9.49"
$ws.Range("E27").Value = "  +2.53%  "
$ws.Range("D28").Value = "33.13"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").Value = "8.58"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").Value = "6.98"
$ws.Range("E30").Value = "  -3.72%  "
$ws.Range("D31").Value = "576.22"
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("D32").Value = "11.05"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").Value = "3.755.30"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "57.37"
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("D38").Value = "35.56"
$ws.Range("E38").Value = "  +9.17%  "
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("E40").Value = "  +2.56%  "
$ws.Range("D41").Value = "2.64"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Value = "0.0₃0677"
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").Value = "3.31"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("E49").Value = "  -1.43%  "
$ws.Range("D50").Value = "131.01"
$ws.Range("E50").Value = "  +2.08%  "
$ws.Range("D51").Value = "2.63"
$ws.Range("E51").Value = "  -1.56%  "

# Restore the default (unstyled) look so the cells match the rest of
# the data cells once more.
$priceRange.Style = "Normal"

